$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Beachwood Golf Course")

# ---------------------------------------------------------------------------
# 1. Apply formatting to the new rows (16-22) by copying styles from the
#    equivalent "template" cells already used elsewhere on this sheet, so the
#    same style records (bold, date, time, etc.) are reused.
# ---------------------------------------------------------------------------

# Style "bold label" (same as A2) -> blank separator row 16, plus most label
# cells across the new block.
$ws.Range("A2").Copy()
$ws.Range("A16:L16").PasteSpecial(-4122)
$ws.Range("N16:W16").PasteSpecial(-4122)
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("B18:C18").PasteSpecial(-4122)
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("M18").PasteSpecial(-4122)
$ws.Range("W18").PasteSpecial(-4122)
$ws.Range("X18").PasteSpecial(-4122)
$ws.Range("M19").PasteSpecial(-4122)
$ws.Range("W19").PasteSpecial(-4122)
$ws.Range("X19").PasteSpecial(-4122)
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("K20").PasteSpecial(-4122)
$ws.Range("U20:V20").PasteSpecial(-4122)
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A22").PasteSpecial(-4122)

# Style used for M4/M9 (bold total w/ fill+border flags)
$ws.Range("M4").Copy()
$ws.Range("M16").PasteSpecial(-4122)

# Style used for K11/L11/U11/V11 (explicit default font on score cells)
$ws.Range("K11").Copy()
$ws.Range("K18").PasteSpecial(-4122)
$ws.Range("L18").PasteSpecial(-4122)
$ws.Range("U18").PasteSpecial(-4122)
$ws.Range("V18").PasteSpecial(-4122)
$ws.Range("K19").PasteSpecial(-4122)
$ws.Range("L19").PasteSpecial(-4122)
$ws.Range("U19").PasteSpecial(-4122)
$ws.Range("V19").PasteSpecial(-4122)

# Style used for the round date (A10)
$ws.Range("A10").Copy()
$ws.Range("A17").PasteSpecial(-4122)

# Style used for the round start/end time (B10:C10)
$ws.Range("B10").Copy()
$ws.Range("B17:C17").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Row 16 - blank separator row (formatting only, no values)
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 3. Row 17 - round date / start-time / end-time
# ---------------------------------------------------------------------------
$ws.Range("A17").Value = 44635
$ws.Range("B17").Value = 0.38541666666666669
$ws.Range("C17").Value = 0.58333333333333337

# ---------------------------------------------------------------------------
# 4. Row 18 - "White" tee scores
# ---------------------------------------------------------------------------
$ws.Range("A18").Value = "White"
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = 5
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 6
$ws.Range("I18").Value = 5
$ws.Range("J18").Value = 5
$ws.Range("K18").Value = 6
$ws.Range("L18").Value = 9
$ws.Range("M18").Formula = "=SUM(D18:L18)"
$ws.Range("N18").Value = 7
$ws.Range("O18").Value = 6
$ws.Range("P18").Value = 5
$ws.Range("Q18").Value = 5
$ws.Range("R18").Value = 5
$ws.Range("S18").Value = 6
$ws.Range("T18").Value = 5
$ws.Range("U18").Value = 5
$ws.Range("V18").Value = 5
$ws.Range("W18").Formula = "=SUM(N18:V18)"
$ws.Range("X18").Formula = "=W18+M18"

# ---------------------------------------------------------------------------
# 5. Row 19 - putts
# ---------------------------------------------------------------------------
$ws.Range("A19").Value = "putts"
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 2
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 2
$ws.Range("I19").Value = 2
$ws.Range("J19").Value = 2
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 1
$ws.Range("M19").Formula = "=SUM(D19:L19)"
$ws.Range("N19").Value = 3
$ws.Range("O19").Value = 2
$ws.Range("P19").Value = 2
$ws.Range("Q19").Value = 2
$ws.Range("R19").Value = 1
$ws.Range("S19").Value = 2
$ws.Range("T19").Value = 2
$ws.Range("U19").Value = 1
$ws.Range("V19").Value = 2
$ws.Range("W19").Formula = "=SUM(N19:V19)"
$ws.Range("X19").Formula = "=W19+M19"

# ---------------------------------------------------------------------------
# 6. Row 20 - penalties
# ---------------------------------------------------------------------------
$ws.Range("A20").Value = "penalties"
$ws.Range("P20").Value = " "
$ws.Range("S20").Value = "W"

# ---------------------------------------------------------------------------
# 7. Row 21 - Fairway hit
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = "Fairway"
$ws.Range("D21").Value = "N"
$ws.Range("E21").Value = "Y"
$ws.Range("F21").Value = "Y"
$ws.Range("H21").Value = "N"
$ws.Range("I21").Value = "Y"
$ws.Range("J21").Value = "N"
$ws.Range("L21").Value = "N"
$ws.Range("N21").Value = "N"
$ws.Range("O21").Value = "N"
$ws.Range("P21").Value = "N"
$ws.Range("Q21").Value = "N"
$ws.Range("R21").Value = "Y"
$ws.Range("T21").Value = "Y"
$ws.Range("U21").Value = "N"
$ws.Range("V21").Value = "N"

# ---------------------------------------------------------------------------
# 8. Row 22 - "80 <" strokes-inside-80 tally
# ---------------------------------------------------------------------------
$ws.Range("A22").Value = "80 <"
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 4
$ws.Range("I22").Value = 3
$ws.Range("J22").Value = 3
$ws.Range("K22").Value = 5
$ws.Range("L22").Value = 5
$ws.Range("M22").Formula = "=SUM(D22:L22)"
$ws.Range("N22").Value = 5
$ws.Range("O22").Value = 4
$ws.Range("P22").Value = 3
$ws.Range("Q22").Value = 3
$ws.Range("R22").Value = 3
$ws.Range("S22").Value = 4
$ws.Range("T22").Value = 3
$ws.Range("U22").Value = 2
$ws.Range("V22").Value = 3
$ws.Range("W22").Formula = "=SUM(N22:V22)"
$ws.Range("X22").Formula = "=M22+W22"

# ---------------------------------------------------------------------------
# 9. Selection / active sheet housekeeping to match the new layout.
#    Activating "Beachwood Golf Course" makes it tabSelected and clears the
#    flag on whichever sheet previously had it (Lochmere).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("W22").Select() | Out-Null
